# Fixed naive component forecaster bug - Presentation state 11.02.
# This script updates the error-series matrix values (columns B:K, rows 2:16)
# on the active worksheet to reflect the corrected naive QoQ forecaster output.
# The corrected data extends the staircase of available forecast horizons by
# one extra cell per row starting at row 7 (new cells K7, J8, I9, H10, G11,
# F12, E13, D14, C15, B16 are populated that were previously empty).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.08873594589893813
$ws.Cells.Item(2, 3).Value = 0.5038494199792003
$ws.Cells.Item(2, 4).Value = -0.0804002763435806
$ws.Cells.Item(2, 5).Value = 0.7915875151041314
$ws.Cells.Item(2, 6).Value = 0.7713506143760337
$ws.Cells.Item(2, 7).Value = 0.3749201162859844
$ws.Cells.Item(2, 8).Value = 0.5168869453501342
$ws.Cells.Item(2, 9).Value = 0.6892623450793038
$ws.Cells.Item(2, 10).Value = 0.1844339823288103
$ws.Cells.Item(2, 11).Value = 0.4621062954844631
$ws.Cells.Item(3, 2).Value = 0.526277399612209
$ws.Cells.Item(3, 3).Value = -0.07779673678113191
$ws.Cells.Item(3, 4).Value = 0.7695198294501161
$ws.Cells.Item(3, 5).Value = 0.7653686340716788
$ws.Cells.Item(3, 6).Value = 0.3695768661035733
$ws.Cells.Item(3, 7).Value = 0.506518411979768
$ws.Cells.Item(3, 8).Value = 0.68078502295568
$ws.Cells.Item(3, 9).Value = 0.1766616740219995
$ws.Cells.Item(3, 10).Value = 0.4534823219514945
$ws.Cells.Item(3, 11).Value = 0.2839431369332225
$ws.Cells.Item(4, 2).Value = -0.05637216532391182
$ws.Cells.Item(4, 3).Value = 0.8896976782493284
$ws.Cells.Item(4, 4).Value = 0.6761742896578956
$ws.Cells.Item(4, 5).Value = 0.3398498622549955
$ws.Cells.Item(4, 6).Value = 0.5151445320096781
$ws.Cells.Item(4, 7).Value = 0.6548448495302448
$ws.Cells.Item(4, 8).Value = 0.153883110993772
$ws.Cells.Item(4, 9).Value = 0.4398642868028766
$ws.Cells.Item(4, 10).Value = 0.2654223397480467
$ws.Cells.Item(4, 11).Value = 0.570669944985061
$ws.Cells.Item(5, 2).Value = 0.847377045928939
$ws.Cells.Item(5, 3).Value = 0.6469698158021624
$ws.Cells.Item(5, 4).Value = 0.3368292624500743
$ws.Cells.Item(5, 5).Value = 0.4962832483981977
$ws.Cells.Item(5, 6).Value = 0.6345137184650405
$ws.Cells.Item(5, 7).Value = 0.1393529555595242
$ws.Cells.Item(5, 8).Value = 0.4232232413106087
$ws.Cells.Item(5, 9).Value = 0.2478384943192965
$ws.Cells.Item(5, 10).Value = 0.5541960614550182
$ws.Cells.Item(5, 11).Value = -0.05728328644410208
$ws.Cells.Item(6, 2).Value = 0.9871675564200725
$ws.Cells.Item(6, 3).Value = 0.4122003242340114
$ws.Cells.Item(6, 4).Value = 0.3051899620851986
$ws.Cells.Item(6, 5).Value = 0.6592076310517737
$ws.Cells.Item(6, 6).Value = 0.1487323591158202
$ws.Cells.Item(6, 7).Value = 0.3684124426992176
$ws.Cells.Item(6, 8).Value = 0.2289444034306267
$ws.Cells.Item(6, 9).Value = 0.5393320606399725
$ws.Cells.Item(6, 10).Value = -0.0859949970734728
$ws.Cells.Item(6, 11).Value = 0.6071339948549791
$ws.Cells.Item(7, 2).Value = 0.862895196224262
$ws.Cells.Item(7, 3).Value = 0.352300664297557
$ws.Cells.Item(7, 4).Value = 0.4189244002609654
$ws.Cells.Item(7, 5).Value = 0.184677440181683
$ws.Cells.Item(7, 6).Value = 0.3663616852596248
$ws.Cells.Item(7, 7).Value = 0.1523761639945965
$ws.Cells.Item(7, 8).Value = 0.511932666264689
$ws.Cells.Item(7, 9).Value = -0.1122840472711982
$ws.Cells.Item(7, 10).Value = 0.5637367041416466
$ws.Cells.Item(7, 11).Value = 0.2970525035592049
$ws.Cells.Item(8, 2).Value = 0.6646262512210954
$ws.Cells.Item(8, 3).Value = 0.5522131399964898
$ws.Cells.Item(8, 4).Value = 0.005599018365491398
$ws.Cells.Item(8, 5).Value = 0.394555924030192
$ws.Cells.Item(8, 6).Value = 0.18840147518699
$ws.Cells.Item(8, 7).Value = 0.4743910960604755
$ws.Cells.Item(8, 8).Value = -0.1178263863585594
$ws.Cells.Item(8, 9).Value = 0.5679120330803951
$ws.Cells.Item(8, 10).Value = 0.2858677898194339
$ws.Cells.Item(9, 2).Value = 0.787803631104331
$ws.Cells.Item(9, 3).Value = 0.09027775923980097
$ws.Cells.Item(9, 4).Value = 0.2489555573964748
$ws.Cells.Item(9, 5).Value = 0.1983552180462326
$ws.Cells.Item(9, 6).Value = 0.4897562657600204
$ws.Cells.Item(9, 7).Value = -0.1566747213159825
$ws.Cells.Item(9, 8).Value = 0.5520662240532093
$ws.Cells.Item(9, 9).Value = 0.2775335613519331
$ws.Cells.Item(10, 2).Value = 0.4013019457211838
$ws.Cells.Item(10, 3).Value = 0.36604433180767
$ws.Cells.Item(10, 4).Value = 0.03589107659666579
$ws.Cells.Item(10, 5).Value = 0.518715216225222
$ws.Cells.Item(10, 6).Value = -0.1209318488610789
$ws.Cells.Item(10, 7).Value = 0.5207385776695821
$ws.Cells.Item(10, 8).Value = 0.2743085116504074
$ws.Cells.Item(11, 2).Value = 0.6128695092117844
$ws.Cells.Item(11, 3).Value = 0.0535469441345553
$ws.Cells.Item(11, 4).Value = 0.4240929771142275
$ws.Cells.Item(11, 5).Value = -0.0887144606125988
$ws.Cells.Item(11, 6).Value = 0.5331267034972994
$ws.Cells.Item(11, 7).Value = 0.2534447081011285
$ws.Cells.Item(12, 2).Value = 0.2932233035507672
$ws.Cells.Item(12, 3).Value = 0.509117220583441
$ws.Cells.Item(12, 4).Value = -0.2051990389706129
$ws.Cells.Item(12, 5).Value = 0.5482319751491519
$ws.Cells.Item(12, 6).Value = 0.2766837437271186
$ws.Cells.Item(13, 2).Value = 0.6739772976175282
$ws.Cells.Item(13, 3).Value = -0.1916617667226967
$ws.Cells.Item(13, 4).Value = 0.4852590561591889
$ws.Cells.Item(13, 5).Value = 0.2867219094086165
$ws.Cells.Item(14, 2).Value = 0.06218727514271133
$ws.Cells.Item(14, 3).Value = 0.5845771063412253
$ws.Cells.Item(14, 4).Value = 0.1751453671933744
$ws.Cells.Item(15, 2).Value = 0.6286367975806744
$ws.Cells.Item(15, 3).Value = 0.1965658720679752
$ws.Cells.Item(16, 2).Value = 0.4328090033804217
